# Add a fourth test-case block (TC004) to the CreateEmployee sheet, mirroring
# the existing TC001/TC002/TC003 blocks: a header row, a name/job header row,
# and one data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: section header, styled like the existing "TC00x_..." header cells (A1/A6/A12)
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = "TC004__ValidateCreateAPILatency"

# Row 18: "name"/"job" column header, styled like the existing header rows (A2:B2/A7:B7/A13:B13)
$ws.Range("A13:B13").Copy() | Out-Null
$ws.Range("A18:B18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = "name"
$ws.Range("B18").Value = "job"

# Row 19: data row, styled like the existing data rows (A14:B14)
$ws.Range("A14:B14").Copy() | Out-Null
$ws.Range("A19:B19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = "Popy"
$ws.Range("B19").Value = "PM"

$ws.Application.CutCopyMode = $false

# Update the active selection to reflect where the new data was entered.
$ws.Range("A17").Select() | Out-Null
